$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row above the old row 4 ("Number of disability
#    persons") so that we end up with two data rows instead of one:
#    row4 -> "family with disabilities Persons "
#    row5 -> "disabilities Persons " (used to be row4)
#    row6 -> Source note (used to be row5)
# ---------------------------------------------------------------------
$ws.Rows(4).Insert()

# ---------------------------------------------------------------------
# 2. Title row (row 1): new wording, merged across A1:I1, bold, wrapped,
#    centered both ways, taller row.
# ---------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Abasha Municipality"
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").WrapText = $true
$ws.Rows(1).RowHeight = 51

# ---------------------------------------------------------------------
# 3. Row 2 ("(End of year, persons)") keeps its text, but reverts to the
#    default row height (no explicit custom height any more).
# ---------------------------------------------------------------------
$ws.Rows(2).RowHeight = 14.5

# ---------------------------------------------------------------------
# 4. Row 4 - new "family with disabilities Persons " row.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(8).Weight = 2

$rowVals4 = @(622,605,551,556,556,565,571,583)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $rowVals4[$i]
    $cell.NumberFormat = "#\ ##0"
}
$ws.Rows(4).RowHeight = 24.75

# ---------------------------------------------------------------------
# 5. Row 5 - "disabilities Persons " row (was row4 before the insert).
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2

$rowVals5 = @(694,683,626,632,629,633,638,654)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $rowVals5[$i]
    $cell.NumberFormat = "#\ ##0"
}
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
$ws.Rows(5).RowHeight = 21

# ---------------------------------------------------------------------
# 6. Row 6 - Source note row (was row5). Text is unchanged, merge range
#    shifts down to A6:H6.
# ---------------------------------------------------------------------
$ws.Rows(6).RowHeight = 27.75

# ---------------------------------------------------------------------
# 7. Column widths - only column A should keep a custom width now.
# ---------------------------------------------------------------------
$ws.Columns("B:R").ColumnWidth = 8.43
$ws.Columns("A").ColumnWidth = 20.83

# ---------------------------------------------------------------------
# 8. Selection as saved in the file.
# ---------------------------------------------------------------------
$ws.Range("A1:I1").Select()

Write-Output "edit complete"
